$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 941.9091
$ws.Range("I33").Value = 47.375
$ws.Range("J33").Value = 3327.3333
$ws.Range("K33").Value = 47.375
$ws.Range("L33").Value = 3327.3333
$ws.Range("M33").Value = 181.625
$ws.Range("N33").Value = -3785.3333
# Row 51
$ws.Range("H51").Value = 9443.666999999999
$ws.Range("I51").Value = 7500
$ws.Range("J51").Value = 9999
$ws.Range("K51").Value = 7500
$ws.Range("L51").Value = 9999
$ws.Range("M51").Value = -7016
$ws.Range("N51").Value = -10967
# Row 69
$ws.Range("H69").Value = 5000
$ws.Range("J69").Value = 5000
$ws.Range("L69").Value = 15000
$ws.Range("N69").Value = -16748
# Row 72
$ws.Range("H72").Value = 5000
$ws.Range("J72").Value = 5000
$ws.Range("L72").Value = 45000
$ws.Range("N72").Value = -53736
# Row 80
$ws.Range("H80").Value = 1833.0714
$ws.Range("I80").Value = 1180.1666
$ws.Range("K80").Value = 3540.4998
$ws.Range("M80").Value = -2542.4998
# Row 83
$ws.Range("H83").Value = 1833.0714
$ws.Range("I83").Value = 1180.1666
$ws.Range("K83").Value = 10621.4994
$ws.Range("M83").Value = -5629.499400000001
# Row 138
$ws.Range("H138").Value = 1615.6666
$ws.Range("J138").Value = 3000
$ws.Range("L138").Value = 9000
$ws.Range("N138").Value = -19280

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
# Row 5
$ws.Range("H5").Value = 127.8
$ws.Range("J5").Value = 104.5
$ws.Range("L5").Value = 104.5
$ws.Range("N5").Value = -328.5
# Row 32
$ws.Range("H32").Value = 5633.1113
$ws.Range("I32").Value = 5537.25
$ws.Range("J32").Value = 6400
$ws.Range("K32").Value = 5537.25
$ws.Range("L32").Value = 6400
$ws.Range("M32").Value = -5250.25
$ws.Range("N32").Value = -6974
# Row 95
$ws.Range("H95").Value = 16499.5
$ws.Range("J95").Value = 16499.5
$ws.Range("L95").Value = 16499.5
$ws.Range("N95").Value = -21991.5
# Row 97
$ws.Range("H97").Value = 784.6667
$ws.Range("I97").Value = 794
$ws.Range("K97").Value = 794
$ws.Range("M97").Value = -298
# Row 102
$ws.Range("H102").Value = 1983.5
$ws.Range("I102").Value = 1869.3334
$ws.Range("K102").Value = 1869.3334
$ws.Range("M102").Value = -247.3334
# Row 130
$ws.Range("H130").Value = 42629.332
$ws.Range("J130").Value = 42629.332
$ws.Range("L130").Value = 42629.332
$ws.Range("N130").Value = -52669.332

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 127.8
$ws.Range("J4").Value = 104.5
$ws.Range("L4").Value = 104.5
$ws.Range("N4").Value = -334.5
# Row 20
$ws.Range("H20").Value = 1863.4286
$ws.Range("I20").Value = 1441.5
$ws.Range("J20").Value = 4395
$ws.Range("K20").Value = 1441.5
$ws.Range("L20").Value = 4395
$ws.Range("M20").Value = -1194.5
$ws.Range("N20").Value = -4889
# Row 99
$ws.Range("H99").Value = 2378.3333
$ws.Range("J99").Value = 2000
$ws.Range("L99").Value = 2000
$ws.Range("N99").Value = -4996
# Row 124
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()
# Row 125
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 294.2
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()
# Row 22
$ws.Range("H22").Value = 999.3333
$ws.Range("I22").Value = 999.5
$ws.Range("J22").Value = 999
$ws.Range("K22").Value = 999.5
$ws.Range("L22").Value = 999
$ws.Range("M22").Value = -649.5
$ws.Range("N22").Value = -1699
# Row 57
$ws.Range("H57").Value = 6000
$ws.Range("J57").Value = 6000
$ws.Range("L57").Value = 6000
$ws.Range("N57").Value = -7120
# Row 74
$ws.Range("H74").Value = 23666.666
$ws.Range("J74").Value = 32500
$ws.Range("L74").Value = 32500
$ws.Range("N74").Value = -34248
# Row 75
$ws.Range("H75").Value = 44444
$ws.Range("J75").Value = 44444
$ws.Range("L75").Value = 44444
$ws.Range("N75").Value = -46440
# Row 77
$ws.Range("H77").Value = 23666.666
$ws.Range("J77").Value = 32500
$ws.Range("L77").Value = 97500
$ws.Range("N77").Value = -106236
# Row 78
$ws.Range("H78").Value = 44444
$ws.Range("J78").Value = 44444
$ws.Range("L78").Value = 133332
$ws.Range("N78").Value = -143316
# Row 80
$ws.Range("H80").Value = 29999
$ws.Range("J80").Value = 29999
$ws.Range("L80").Value = 29999
$ws.Range("N80").Value = -32245
# Row 83
$ws.Range("H83").Value = 29999
$ws.Range("J83").Value = 29999
$ws.Range("L83").Value = 89997
$ws.Range("N83").Value = -101229
# Row 94
$ws.Range("H94").Value = 7146.5
$ws.Range("I94").Value = 7146.5
$ws.Range("K94").Value = 7146.5
$ws.Range("M94").Value = -6695.5
# Row 99
$ws.Range("H99").Value = 4983
$ws.Range("I99").Value = 4237
$ws.Range("J99").Value = 6475
$ws.Range("K99").Value = 4237
$ws.Range("L99").Value = 6475
$ws.Range("M99").Value = -2739
$ws.Range("N99").Value = -9471
# Row 126
$ws.Range("H126").Value = 4983
$ws.Range("I126").Value = 4237
$ws.Range("J126").Value = 6475
$ws.Range("K126").Value = 12711
$ws.Range("L126").Value = 19425
$ws.Range("M126").Value = -10241
$ws.Range("N126").Value = -24365
# Row 129
$ws.Range("H129").Value = 0
$ws.Range("I129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("M129").ClearContents()
# Row 130
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 22.777779
$ws.Range("J2").Value = 17.5
$ws.Range("L2").Value = 105
$ws.Range("N2").Value = -331
# Row 4
$ws.Range("H4").Value = 102
$ws.Range("I4").Value = 102
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 306
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -194
$ws.Range("N4").ClearContents()
# Row 17
$ws.Range("H17").Value = 64.61539
$ws.Range("I17").Value = 52.1
$ws.Range("K17").Value = 156.3
$ws.Range("M17").Value = 12.69999999999999
# Row 132
$ws.Range("H132").Value = 1698.3334
$ws.Range("I132").Value = 1698.3334
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 15285.0006
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -12755.0006
$ws.Range("N132").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 63
$ws.Range("H63").Value = 10000
$ws.Range("J63").Value = 10000
$ws.Range("L63").Value = 10000
$ws.Range("N63").Value = -11372
# Row 66
$ws.Range("H66").Value = 10000
$ws.Range("J66").Value = 10000
$ws.Range("L66").Value = 30000
$ws.Range("N66").Value = -36864
# Row 107
$ws.Range("H107").Value = 300.5
$ws.Range("I107").Value = 400
$ws.Range("K107").Value = 400
$ws.Range("M107").Value = 1520
# Row 113
$ws.Range("H113").Value = 4955.4443
$ws.Range("I113").Value = 4583.7144
$ws.Range("J113").Value = 6256.5
$ws.Range("K113").Value = 4583.7144
$ws.Range("L113").Value = 6256.5
$ws.Range("M113").Value = -2413.7144
$ws.Range("N113").Value = -10596.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 1265.8
$ws.Range("I7").Value = 1265.8
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1265.8
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -1153.8
$ws.Range("N7").ClearContents()
# Row 126
$ws.Range("H126").Value = 1265.8
$ws.Range("I126").Value = 1265.8
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 3797.4
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -1327.4
$ws.Range("N126").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 2
$ws.Range("H2").Value = 90000
$ws.Range("I2").Value = 90000
$ws.Range("K2").Value = 90000
$ws.Range("M2").Value = -89888
# Row 62
$ws.Range("H62").Value = 1916.3334
$ws.Range("I62").Value = 2124.5
$ws.Range("J62").Value = 1500
$ws.Range("K62").Value = 2124.5
$ws.Range("L62").Value = 1500
$ws.Range("M62").Value = -1500.5
$ws.Range("N62").Value = -2748
# Row 65
$ws.Range("H65").Value = 1916.3334
$ws.Range("I65").Value = 2124.5
$ws.Range("J65").Value = 1500
$ws.Range("K65").Value = 10622.5
$ws.Range("L65").Value = 7500
$ws.Range("M65").Value = -7502.5
$ws.Range("N65").Value = -13740
# Row 107
$ws.Range("H107").Value = 800
$ws.Range("I107").Value = 800
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2400
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -480
$ws.Range("N107").ClearContents()
# Row 130
$ws.Range("H130").Value = 60000
$ws.Range("J130").Value = 60000
$ws.Range("L130").Value = 60000
$ws.Range("N130").Value = -70040
# Row 136
$ws.Range("H136").Value = 2431.6667
$ws.Range("I136").Value = 1743.7273
$ws.Range("J136").Value = 9999
$ws.Range("K136").Value = 5231.1819
$ws.Range("L136").Value = 29997
$ws.Range("M136").Value = -2681.1819
$ws.Range("N136").Value = -35097
